$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string need Text format
# so Excel stores them as text (matching the source data) instead of coercing to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.070.79'
$ws.Range("E2").Value = '  -2.84%  '
$ws.Range("D3").Value = '1.844.28'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '0.6919'
$ws.Range("E5").Value = '  -6.89%  '
$ws.Range("D6").Value = '237.36'
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '0.07598'
$ws.Range("E8").Value = '  +5.24%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.3034'
$ws.Range("E9").Value = '  -3.72%  '
$ws.Range("D10").Value = '23.22'
$ws.Range("E10").Value = '  -5.81%  '
$ws.Range("D11").Value = '0.08097'
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").Value = '1.874.97'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '0.7197'
$ws.Range("E13").Value = '  -4.06%  '
$ws.Range("D14").Value = '5.164'
$ws.Range("E14").Value = '  -3.98%  '
$ws.Range("D15").Value = '88.88'
$ws.Range("E15").Value = '  -3.65%  '
$ws.Range("D16").Value = '29.168.00'
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '241.85'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '5.757'
$ws.Range("E18").Value = '  -5.94%  '
$ws.Range("D19").Value = '0.000007688'
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").Value = '13.04'
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '2.103.15'
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '7.589'
$ws.Range("E24").Value = '  -5.22%  '
$ws.Range("D25").Value = '8.969'
$ws.Range("E25").Value = '  -3.48%  '
$ws.Range("D26").Value = '161.00'
$ws.Range("E26").Value = '  -2.71%  '
$ws.Range("D27").Value = '0.1455'
$ws.Range("E27").Value = '  -5.69%  '
$ws.Range("D28").Value = '18.03'
$ws.Range("E28").Value = '  -3.37%  '
$ws.Range("D29").Value = '1.924'
$ws.Range("E29").Value = '  -4.68%  '
$ws.Range("D30").Value = '1.380'
$ws.Range("E30").Value = '  -8.10%  '
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("D32").Value = '1.491'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").Value = '4.042'
$ws.Range("E33").Value = '  -4.19%  '
$ws.Range("D34").Value = '0.05225'
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("D35").Value = '1.181'
$ws.Range("E35").Value = '  -4.32%  '
$ws.Range("D36").Value = '0.7080'
$ws.Range("E36").Value = '  -5.59%  '
$ws.Range("D37").Value = '0.9991'
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("D38").Value = '2.663'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").Value = '0.01851'
$ws.Range("E39").Value = '  -5.76%  '
$ws.Range("D40").Value = '2.687'
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").Value = '0.9140'
$ws.Range("E41").Value = '  +5.77%  '
$ws.Range("D42").Value = '5.949'
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("D43").Value = '0.4273'
$ws.Range("E43").Value = '  -5.42%  '
$ws.Range("D44").Value = '69.64'
$ws.Range("E44").Value = '  -3.85%  '
$ws.Range("D45").Value = '1.044.85'
$ws.Range("E45").Value = '  -6.34%  '
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").Value = '102.31'
$ws.Range("E47").Value = '  -1.92%  '
$ws.Range("D48").Value = '7.184'
$ws.Range("E48").Value = '  -5.74%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.738'
$ws.Range("E49").Value = '  -6.72%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.000.63'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("D51").Value = '9.216'
$ws.Range("E51").Value = '  -3.14%  '
